$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered style)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-52: team record (Wins=55, Losses=107, Ties=0) for every player row
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 55   # AD
    $ws.Cells.Item($r, 31).Value = 107  # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
